$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 81878556 file (row 4)
$wsOverview.Range("G4").Value = "2016-08-19 08:16:58"
$wsOverview.Range("G5").Value = "2016-08-19 08:16:58"

# zh-cn sheet, row 4 (81878556 file)
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-19 08:16:54"
$wsZhCn.Range("K4").Value = "2016-08-19 08:17:15"

# zh-cn sheet, row 5 (cf2bdbd0 file)
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H5").Value = "2016-08-19 08:16:54"
$wsZhCn.Range("K5").Value = "2016-08-19 08:17:15"

# de-de sheet, row 4 (81878556 file)
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-19 08:16:58"
$wsDeDe.Range("K4").Value = "2016-08-19 08:17:21"

# de-de sheet, row 5 (cf2bdbd0 file)
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H5").Value = "2016-08-19 08:16:58"
$wsDeDe.Range("K5").Value = "2016-08-19 08:17:21"
